$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21 : 2020-12-29, 3h, "Raspberry fuer VPN..." ---
$ws.Range("A21").Value = 44194
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)

$ws.Range("B21").Value = 3

$ws.Range("C21").Formula = "=C20+B21"

$ws.Range("D21").Value = "Raspberry für VPN eingerichtet und über wireguard verbunden"
$ws.Range("D20").Copy()
$ws.Range("D21").PasteSpecial(-4122)

$ws.Rows("21").RowHeight = 30

# --- Row 22 : 2020-12-30, 4h, "Versucht VPN verbindung..." ---
$ws.Range("A22").Value = 44195
$ws.Range("A20").Copy()
$ws.Range("A22").PasteSpecial(-4122)

$ws.Range("B22").Value = 4

$ws.Range("C22").Formula = "=C21+B22"

$ws.Range("D22").Value = "Versucht VPN verbindung auf Redpitaya weiterzugeben" + [char]10 + "Entweder fehler beim bridgen oder Redpitaya bekam immer eigene IP-Adresse ohne VPN"
$ws.Range("D20").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Rows("22").RowHeight = 60

$excel.CutCopyMode = $false

# --- Window / selection state ---
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D23").Select()
